# Rename the "department" value for the General English course row
# from "FACULTY OF ENGLISH" to "English".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "English"

# Match the author's final selection position recorded in the saved file.
$ws.Range("G8").Select()
